$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.061.12'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '2.305.33'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.79'
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.04'
$ws.Range('E6').Value = '  -1.16%  '
$ws.Range('E7').Value = '  +4.16%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.69'
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.03'
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '2.663.16'
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').Value = '2.313.87'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('D18').Value = '42.969.62'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.36'
$ws.Range('E19').Value = '  +7.96%  '
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.13'
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.22'
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.58'
$ws.Range('E23').Value = '  +1.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.19'
$ws.Range('E24').Value = '  -1.09%  '
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.71'
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '168.37'
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.15'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.04'
$ws.Range('E30').Value = '  -0.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.92'
$ws.Range('E31').Value = '  -4.18%  '
$ws.Range('E32').Value = '  +4.34%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.80'
$ws.Range('E34').Value = '  +4.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.12'
$ws.Range('E35').Value = '  +3.49%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0690'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.103'
$ws.Range('E38').Value = '  +1.24%  '
$ws.Range('E39').Value = '  +1.16%  '
$ws.Range('E40').Value = '  +2.30%  '
$ws.Range('E41').Value = '  -2.45%  '
$ws.Range('D42').Value = '2.012.56'
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.16'
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.14'
$ws.Range('E45').Value = '  -3.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.46'
$ws.Range('E46').Value = '  +0.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.83'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.49'
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('D49').Value = '2.528.90'
$ws.Range('E49').Value = '  +0.60%  '
$ws.Range('E50').Value = '  +0.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.43'
$ws.Range('E51').Value = '  +4.87%  '
